$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the comment text in H5 (shared string referenced by this cell)
$ws.Range("H5").Value = "Autonomous 27-minute file matched in time to TriOS FICE22 stations. Min/max rotator: -126/+42"

# Update row 5 height to accommodate the longer wrapped text
$ws.Rows.Item(5).RowHeight = 68
